$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AlphaFiberF-HW45")

# --- Update shared-string driven labels / rename "Thomas Hex" -> "Matthies Hex" is folded
# into the full data rewrite below (B-column values), since every row shifts down by two
# to make room for the two new pole-figure types ("Holden" and "Rizzie Spiral") that were
# added after "Spiral5" and the simulation was rerun with updated results.

$data = New-Object 'object[,]' 30,20
$data[0,0] = 0
$data[0,1] = 'HKL'
$data[0,2] = '[3, 2, 1]'
$data[0,3] = '[1, 1, 0]'
$data[0,4] = '[3, 1, 0]'
$data[0,5] = '[2, 2, 2]'
$data[0,6] = '[2, 0, 0]'
$data[0,7] = '[2, 2, 0]'
$data[0,8] = '[2, 1, 1]'
$data[0,9] = '[4, 0, 0]'
$data[0,10] = '1Pair-A'
$data[0,11] = '1Pair-B'
$data[0,12] = '2Pairs-A'
$data[0,13] = '2Pairs-B'
$data[0,14] = '3Pairs-A'
$data[0,15] = '3Pairs-B'
$data[0,16] = '3Pairs-C'
$data[0,17] = '4Pairs'
$data[0,18] = '5A4F'
$data[0,19] = 'MaxUnique'
$data[1,0] = 1
$data[1,1] = 'Spiral5'
$data[1,2] = 1.00000039235269
$data[1,3] = 1.00000039235269
$data[1,4] = 0.9999991394554139
$data[1,5] = 1.000001045778407
$data[1,6] = 0.9999984326946321
$data[1,7] = 1.00000039235269
$data[1,8] = 1.00000039235269
$data[1,9] = 0.9999984326946321
$data[1,10] = 1.00000039235269
$data[1,11] = 1.00000039235269
$data[1,12] = 0.9999994125236611
$data[1,13] = 0.9999994125236611
$data[1,14] = 0.999999321500912
$data[1,15] = 0.9999997391333375
$data[1,16] = 0.9999997391333375
$data[1,17] = 0.9999999024381756
$data[1,18] = 0.9999999024381756
$data[1,19] = 0.9999999658310871
$data[2,0] = 2
$data[2,1] = 'Holden'
$data[2,2] = 0.9998076608020819
$data[2,3] = 0.9998076608020819
$data[2,4] = 1.000423146952147
$data[2,5] = 0.9994870940347622
$data[2,6] = 1.000769358945351
$data[2,7] = 0.9998076608020819
$data[2,8] = 0.9998076608020819
$data[2,9] = 1.000769358945351
$data[2,10] = 0.9998076608020819
$data[2,11] = 0.9998076608020819
$data[2,12] = 1.000288509873716
$data[2,13] = 1.000288509873716
$data[2,14] = 1.00033338889986
$data[2,15] = 1.000128226849838
$data[2,16] = 1.000128226849838
$data[2,17] = 1.000048085337899
$data[2,18] = 1.000048085337899
$data[2,19] = 1.000017097056418
$data[3,0] = 3
$data[3,1] = 'Rizzie Spiral'
$data[3,2] = 1.003182996314844
$data[3,3] = 1.003182996314844
$data[3,4] = 0.9929974386417181
$data[3,5] = 1.008487957807861
$data[3,6] = 0.9872680631307793
$data[3,7] = 1.003182996314844
$data[3,8] = 1.003182996314844
$data[3,9] = 0.9872680631307793
$data[3,10] = 1.003182996314844
$data[3,11] = 1.003182996314844
$data[3,12] = 0.9952255297228119
$data[3,13] = 0.9952255297228119
$data[3,14] = 0.9944828326957805
$data[3,15] = 0.9978780185868228
$data[3,16] = 0.9978780185868228
$data[3,17] = 0.9992042630188283
$data[3,18] = 0.9992042630188283
$data[3,19] = 0.9997170747541486
$data[4,0] = 4
$data[4,1] = 'RotRing OmegaMax-90'
$data[4,2] = 1.001083911538769
$data[4,3] = 1.001083911538769
$data[4,4] = 0.9976153918073414
$data[4,5] = 1.002890437022457
$data[4,6] = 0.9956643448017533
$data[4,7] = 1.001083911538769
$data[4,8] = 1.001083911538769
$data[4,9] = 0.9956643448017533
$data[4,10] = 1.001083911538769
$data[4,11] = 1.001083911538769
$data[4,12] = 0.9983741281702611
$data[4,13] = 0.9983741281702611
$data[4,14] = 0.9981212160492877
$data[4,15] = 0.999277389293097
$data[4,16] = 0.999277389293097
$data[4,17] = 0.9997290198545149
$data[4,18] = 0.9997290198545149
$data[4,19] = 0.9999036513746432
$data[5,0] = 5
$data[5,1] = 'Equal Angle'
$data[5,2] = 0.999834895158502
$data[5,3] = 0.999834895158502
$data[5,4] = 1.000363230511528
$data[5,5] = 0.9995597187463992
$data[5,6] = 1.000660421051875
$data[5,7] = 0.999834895158502
$data[5,8] = 0.999834895158502
$data[5,9] = 1.000660421051875
$data[5,10] = 0.999834895158502
$data[5,11] = 0.999834895158502
$data[5,12] = 1.000247658105188
$data[5,13] = 1.000247658105188
$data[5,14] = 1.000286182240635
$data[5,15] = 1.000110070456293
$data[5,16] = 1.000110070456293
$data[5,17] = 1.000041276631845
$data[5,18] = 1.000041276631845
$data[5,19] = 1.000014675964218
$data[6,0] = 6
$data[6,1] = 'Tilt Rotate'
$data[6,2] = 0.9993040803342889
$data[6,3] = 0.9993040803342889
$data[6,4] = 1.001531005075255
$data[6,5] = 0.9981442074371726
$data[6,6] = 1.002783669019255
$data[6,7] = 0.9993040803342889
$data[6,8] = 0.9993040803342889
$data[6,9] = 1.002783669019255
$data[6,10] = 0.9993040803342889
$data[6,11] = 0.9993040803342889
$data[6,12] = 1.001043874676772
$data[6,13] = 1.001043874676772
$data[6,14] = 1.001206251476266
$data[6,15] = 1.000463943229277
$data[6,16] = 1.000463943229277
$data[6,17] = 1.00017397750553
$data[6,18] = 1.00017397750553
$data[6,19] = 1.000061853755758
$data[7,0] = 7
$data[7,1] = 'CLR'
$data[7,2] = 0.9999912783450363
$data[7,3] = 0.9999912783450363
$data[7,4] = 1.000019190950489
$data[7,5] = 0.9999767419877152
$data[7,6] = 1.000034888234707
$data[7,7] = 0.9999912783450363
$data[7,8] = 0.9999912783450363
$data[7,9] = 1.000034888234707
$data[7,10] = 0.9999912783450363
$data[7,11] = 0.9999912783450363
$data[7,12] = 1.000013083289872
$data[7,13] = 1.000013083289872
$data[7,14] = 1.000015119176744
$data[7,15] = 1.000005814974927
$data[7,16] = 1.000005814974927
$data[7,17] = 1.000002180817454
$data[7,18] = 1.000002180817454
$data[7,19] = 1.00000077603467
$data[8,0] = 8
$data[8,1] = 'Rizzie Hex'
$data[8,2] = 0.9999982671232656
$data[8,3] = 0.9999982671232656
$data[8,4] = 1.000003815553753
$data[8,5] = 0.9999953776582707
$data[8,6] = 1.000006934838652
$data[8,7] = 0.9999982671232656
$data[8,8] = 0.9999982671232656
$data[8,9] = 1.000006934838652
$data[8,10] = 0.9999982671232656
$data[8,11] = 0.9999982671232656
$data[8,12] = 1.000002600980959
$data[8,13] = 1.000002600980959
$data[8,14] = 1.000003005838557
$data[8,15] = 1.000001156361728
$data[8,16] = 1.000001156361728
$data[8,17] = 1.000000434052112
$data[8,18] = 1.000000434052112
$data[8,19] = 1.000000154903412
$data[9,0] = 9
$data[9,1] = 'Matthies Hex'
$data[9,2] = 0.9999617214304334
$data[9,3] = 0.9999617214304334
$data[9,4] = 1.0000842147005
$data[9,5] = 0.9998979232219442
$data[9,6] = 1.000153118317757
$data[9,7] = 0.9999617214304334
$data[9,8] = 0.9999617214304334
$data[9,9] = 1.000153118317757
$data[9,10] = 0.9999617214304334
$data[9,11] = 0.9999617214304334
$data[9,12] = 1.000057419874095
$data[9,13] = 1.000057419874095
$data[9,14] = 1.000066351482897
$data[9,15] = 1.000025520392875
$data[9,16] = 1.000025520392875
$data[9,17] = 1.000009570652264
$data[9,18] = 1.000009570652264
$data[9,19] = 1.000003403421917
$data[10,0] = 10
$data[10,1] = 'Tilt Rotate_Partial'
$data[10,2] = 0.9992936413809913
$data[10,3] = 0.9992936413809913
$data[10,4] = 1.001553971795739
$data[10,5] = 0.9981163681590556
$data[10,6] = 1.002825428719375
$data[10,7] = 0.9992936413809913
$data[10,8] = 0.9992936413809913
$data[10,9] = 1.002825428719375
$data[10,10] = 0.9992936413809913
$data[10,11] = 0.9992936413809913
$data[10,12] = 1.001059535050183
$data[10,13] = 1.001059535050183
$data[10,14] = 1.001224347298702
$data[10,15] = 1.000470903827119
$data[10,16] = 1.000470903827119
$data[10,17] = 1.000176588215587
$data[10,18] = 1.000176588215587
$data[10,19] = 1.000062782136191
$data[11,0] = 11
$data[11,1] = 'RotRing OmegaMax-60'
$data[11,2] = 1.000809140556624
$data[11,3] = 1.000809140556624
$data[11,4] = 0.9982198895632023
$data[11,5] = 1.002157714085841
$data[11,6] = 0.9967634287253719
$data[11,7] = 1.000809140556624
$data[11,8] = 1.000809140556624
$data[11,9] = 0.9967634287253719
$data[11,10] = 1.000809140556624
$data[11,11] = 1.000809140556624
$data[11,12] = 0.9987862846409977
$data[11,13] = 0.9987862846409977
$data[11,14] = 0.9985974862817325
$data[11,15] = 0.9994605699462064
$data[11,16] = 0.9994605699462064
$data[11,17] = 0.9997977125988107
$data[11,18] = 0.9997977125988107
$data[11,19] = 0.9999280756740476
$data[12,0] = 12
$data[12,1] = 'Equal Angle_Partial'
$data[12,2] = 0.9998196751157897
$data[12,3] = 0.9998196751157897
$data[12,4] = 1.000396714589474
$data[12,5] = 0.9995191297473706
$data[12,6] = 1.000721305294738
$data[12,7] = 0.9998196751157897
$data[12,8] = 0.9998196751157897
$data[12,9] = 1.000721305294738
$data[12,10] = 0.9998196751157897
$data[12,11] = 0.9998196751157897
$data[12,12] = 1.000270490205264
$data[12,13] = 1.000270490205264
$data[12,14] = 1.000312565
$data[12,15] = 1.000120218508773
$data[12,16] = 1.000120218508773
$data[12,17] = 1.000045082660527
$data[12,18] = 1.000045082660527
$data[12,19] = 1.000016029163158
$data[13,0] = 13
$data[13,1] = 'Rizzie Hex_Partial'
$data[13,2] = 1.00013892584912
$data[13,3] = 1.00013892584912
$data[13,4] = 0.9996943700054491
$data[13,5] = 1.00037046336123
$data[13,6] = 0.9994443108028237
$data[13,7] = 1.00013892584912
$data[13,8] = 1.00013892584912
$data[13,9] = 0.9994443108028237
$data[13,10] = 1.00013892584912
$data[13,11] = 1.00013892584912
$data[13,12] = 0.9997916183259717
$data[13,13] = 0.9997916183259717
$data[13,14] = 0.9997592022191308
$data[13,15] = 0.9999073875003545
$data[13,16] = 0.9999073875003545
$data[13,17] = 0.9999652720875458
$data[13,18] = 0.9999652720875458
$data[13,19] = 0.9999876536194771
$data[14,0] = 14
$data[14,1] = 'ND Single'
$data[14,2] = 0.9987796199999999
$data[14,3] = 0.9987796199999999
$data[14,4] = 1.002684799999999
$data[14,5] = 0.99674564
$data[14,6] = 1.004881500000001
$data[14,7] = 0.9987796199999999
$data[14,8] = 0.9987796199999999
$data[14,9] = 1.004881500000001
$data[14,10] = 0.9987796199999999
$data[14,11] = 0.9987796199999999
$data[14,12] = 1.00183056
$data[14,13] = 1.00183056
$data[14,14] = 1.002115306666666
$data[14,15] = 1.00081358
$data[14,16] = 1.00081358
$data[14,17] = 1.00030509
$data[14,18] = 1.00030509
$data[14,19] = 1.000108466666667
$data[15,0] = 15
$data[15,1] = 'RD Single'
$data[15,2] = 1.0042979
$data[15,3] = 1.0042979
$data[15,4] = 0.99054463
$data[15,5] = 1.0114611
$data[15,6] = 0.98280841
$data[15,7] = 1.0042979
$data[15,8] = 1.0042979
$data[15,9] = 0.98280841
$data[15,10] = 1.0042979
$data[15,11] = 1.0042979
$data[15,12] = 0.9935531550000001
$data[15,13] = 0.9935531550000001
$data[15,14] = 0.9925503133333334
$data[15,15] = 0.9971347366666667
$data[15,16] = 0.9971347366666667
$data[15,17] = 0.9989255275000001
$data[15,18] = 0.9989255275000001
$data[15,19] = 0.9996179733333334
$data[16,0] = 16
$data[16,1] = 'TD Single'
$data[16,2] = 1.007902
$data[16,3] = 1.007902
$data[16,4] = 0.98261561
$data[16,5] = 1.021072
$data[16,6] = 0.9683920100000001
$data[16,7] = 1.007902
$data[16,8] = 1.007902
$data[16,9] = 0.9683920100000001
$data[16,10] = 1.007902
$data[16,11] = 1.007902
$data[16,12] = 0.9881470050000001
$data[16,13] = 0.9881470050000001
$data[16,14] = 0.9863032066666667
$data[16,15] = 0.9947320033333334
$data[16,16] = 0.9947320033333334
$data[16,17] = 0.9980245025000001
$data[16,18] = 0.9980245025000001
$data[16,19] = 0.9992976033333333
$data[17,0] = 17
$data[17,1] = 'Morris Single'
$data[17,2] = 1.0030278
$data[17,3] = 1.0030278
$data[17,4] = 0.9933389299999998
$data[17,5] = 1.008074
$data[17,6] = 0.98788896
$data[17,7] = 1.0030278
$data[17,8] = 1.0030278
$data[17,9] = 0.98788896
$data[17,10] = 1.0030278
$data[17,11] = 1.0030278
$data[17,12] = 0.99545838
$data[17,13] = 0.99545838
$data[17,14] = 0.9947518966666666
$data[17,15] = 0.9979815200000001
$data[17,16] = 0.9979815200000001
$data[17,17] = 0.99924309
$data[17,18] = 0.99924309
$data[17,19] = 0.9997308816666667
$data[18,0] = 18
$data[18,1] = 'Ring Perpendicular to ND'
$data[18,2] = 0.9996074950684928
$data[18,3] = 0.9996074950684928
$data[18,4] = 1.000863501780822
$data[18,5] = 0.9989533334246578
$data[18,6] = 1.001570006164384
$data[18,7] = 0.9996074950684928
$data[18,8] = 0.9996074950684928
$data[18,9] = 1.001570006164384
$data[18,10] = 0.9996074950684928
$data[18,11] = 0.9996074950684928
$data[18,12] = 1.000588750616438
$data[18,13] = 1.000588750616438
$data[18,14] = 1.000680334337899
$data[18,15] = 1.00026166543379
$data[18,16] = 1.00026166543379
$data[18,17] = 1.000098122842466
$data[18,18] = 1.000098122842466
$data[18,19] = 1.000034887762557
$data[19,0] = 19
$data[19,1] = 'Ring Perpendicular to RD'
$data[19,2] = 1.001702715789474
$data[19,3] = 1.001702715789474
$data[19,4] = 0.9962540142105263
$data[19,5] = 1.00454058
$data[19,6] = 0.9931891326315789
$data[19,7] = 1.001702715789474
$data[19,8] = 1.001702715789474
$data[19,9] = 0.9931891326315789
$data[19,10] = 1.001702715789474
$data[19,11] = 1.001702715789474
$data[19,12] = 0.9974459242105264
$data[19,13] = 0.9974459242105264
$data[19,14] = 0.997048620877193
$data[19,15] = 0.9988648547368421
$data[19,16] = 0.9988648547368421
$data[19,17] = 0.99957432
$data[19,18] = 0.99957432
$data[19,19] = 0.9998486457017544
$data[20,0] = 20
$data[20,1] = 'Ring Perpendicular to TD'
$data[20,2] = 1.002888287368421
$data[20,3] = 1.002888287368421
$data[20,4] = 0.9936457931578948
$data[20,5] = 1.007702070526316
$data[20,6] = 0.9884468926315789
$data[20,7] = 1.002888287368421
$data[20,8] = 1.002888287368421
$data[20,9] = 0.9884468926315789
$data[20,10] = 1.002888287368421
$data[20,11] = 1.002888287368421
$data[20,12] = 0.9956675899999998
$data[20,13] = 0.9956675899999998
$data[20,14] = 0.9949936577192982
$data[20,15] = 0.9980744891228067
$data[20,16] = 0.9980744891228067
$data[20,17] = 0.9992779386842103
$data[20,18] = 0.9992779386842103
$data[20,19] = 0.9997432697368419
$data[21,0] = 21
$data[21,1] = 'OffsetFTD'
$data[21,2] = 0.9990088213878943
$data[21,3] = 0.9990088213878943
$data[21,4] = 1.002180598738788
$data[21,5] = 0.9973568568408632
$data[21,6] = 1.003964713376434
$data[21,7] = 0.9990088213878943
$data[21,8] = 0.9990088213878943
$data[21,9] = 1.003964713376434
$data[21,10] = 0.9990088213878943
$data[21,11] = 0.9990088213878943
$data[21,12] = 1.001486767382164
$data[21,13] = 1.001486767382164
$data[21,14] = 1.001718044501039
$data[21,15] = 1.000660785384074
$data[21,16] = 1.000660785384074
$data[21,17] = 1.000247794385029
$data[21,18] = 1.000247794385029
$data[21,19] = 1.000088105519961
$data[22,0] = 22
$data[22,1] = 'OffsetATD'
$data[22,2] = 0.9997296356041152
$data[22,3] = 0.9997296356041152
$data[22,4] = 1.000594802802129
$data[22,5] = 0.9992790312211202
$data[22,6] = 1.001081464766156
$data[22,7] = 0.9997296356041152
$data[22,8] = 0.9997296356041152
$data[22,9] = 1.001081464766156
$data[22,10] = 0.9997296356041152
$data[22,11] = 0.9997296356041152
$data[22,12] = 1.000405550185135
$data[22,13] = 1.000405550185135
$data[22,14] = 1.0004686343908
$data[22,15] = 1.000180245324795
$data[22,16] = 1.000180245324795
$data[22,17] = 1.000067592894625
$data[22,18] = 1.000067592894625
$data[22,19] = 1.000024034266959
$data[23,0] = 23
$data[23,1] = 'OffsetF45'
$data[23,2] = 0.9996038835824849
$data[23,3] = 0.9996038835824849
$data[23,4] = 1.000871462401954
$data[23,5] = 0.9989436836861949
$data[23,6] = 1.00158447801317
$data[23,7] = 0.9996038835824849
$data[23,8] = 0.9996038835824849
$data[23,9] = 1.00158447801317
$data[23,10] = 0.9996038835824849
$data[23,11] = 0.9996038835824849
$data[23,12] = 1.000594180797828
$data[23,13] = 1.000594180797828
$data[23,14] = 1.000686607999203
$data[23,15] = 1.000264081726047
$data[23,16] = 1.000264081726047
$data[23,17] = 1.000099032190156
$data[23,18] = 1.000099032190156
$data[23,19] = 1.000035212474796
$data[24,0] = 24
$data[24,1] = 'OffsetA45'
$data[24,2] = 0.9998964493929572
$data[24,3] = 0.9998964493929572
$data[24,4] = 1.000227815231654
$data[24,5] = 0.9997238617926409
$data[24,6] = 1.000414203344274
$data[24,7] = 0.9998964493929572
$data[24,8] = 0.9998964493929572
$data[24,9] = 1.000414203344274
$data[24,10] = 0.9998964493929572
$data[24,11] = 0.9998964493929572
$data[24,12] = 1.000155326368615
$data[24,13] = 1.000155326368615
$data[24,14] = 1.000179489322962
$data[24,15] = 1.000069034043396
$data[24,16] = 1.000069034043396
$data[24,17] = 1.000025887880786
$data[24,18] = 1.000025887880786
$data[24,19] = 1.000009204757907
$data[25,0] = 25
$data[25,1] = 'OffsetFRD'
$data[25,2] = 1.002423139432898
$data[25,3] = 1.002423139432898
$data[25,4] = 0.9946691126875767
$data[25,5] = 1.006461684153397
$data[25,6] = 0.9903074757539098
$data[25,7] = 1.002423139432898
$data[25,8] = 1.002423139432898
$data[25,9] = 0.9903074757539098
$data[25,10] = 1.002423139432898
$data[25,11] = 1.002423139432898
$data[25,12] = 0.996365307593404
$data[25,13] = 0.996365307593404
$data[25,14] = 0.9957999092914616
$data[25,15] = 0.9983845848732354
$data[25,16] = 0.9983845848732354
$data[25,17] = 0.9993942235131512
$data[25,18] = 0.9993942235131512
$data[25,19] = 0.9997846151489299
$data[26,0] = 26
$data[26,1] = 'OffsetARD'
$data[26,2] = 1.000642642036062
$data[26,3] = 1.000642642036062
$data[26,4] = 0.9985862074418604
$data[26,5] = 1.001713696826491
$data[26,6] = 0.9974294612688168
$data[26,7] = 1.000642642036062
$data[26,8] = 1.000642642036062
$data[26,9] = 0.9974294612688168
$data[26,10] = 1.000642642036062
$data[26,11] = 1.000642642036062
$data[26,12] = 0.9990360516524395
$data[26,13] = 0.9990360516524395
$data[26,14] = 0.9988861035822465
$data[26,15] = 0.9995715817803137
$data[26,16] = 0.9995715817803137
$data[26,17] = 0.9998393468442508
$data[26,18] = 0.9998393468442508
$data[26,19] = 0.9999428819408923
$data[27,0] = 27
$data[27,1] = 'Gaussian Quadrature'
$data[27,2] = 0.9999272577892019
$data[27,3] = 0.9999272577892019
$data[27,4] = 1.000160008983498
$data[27,5] = 0.9998060410711583
$data[27,6] = 1.000290931972496
$data[27,7] = 0.9999272577892019
$data[27,8] = 0.9999272577892019
$data[27,9] = 1.000290931972496
$data[27,10] = 0.9999272577892019
$data[27,11] = 0.9999272577892019
$data[27,12] = 1.000109094880849
$data[27,13] = 1.000109094880849
$data[27,14] = 1.000126066248399
$data[27,15] = 1.000048482516966
$data[27,16] = 1.000048482516966
$data[27,17] = 1.000018176335025
$data[27,18] = 1.000018176335025
$data[27,19] = 1.00000645923246
$data[28,0] = 28
$data[28,1] = 'Michael-CCHex'
$data[28,2] = 1.0004716949089
$data[28,3] = 1.0004716949089
$data[28,4] = 0.9989622680710537
$data[28,5] = 1.001257856116698
$data[28,6] = 0.998113224824261
$data[28,7] = 1.0004716949089
$data[28,8] = 1.0004716949089
$data[28,9] = 0.998113224824261
$data[28,10] = 1.0004716949089
$data[28,11] = 1.0004716949089
$data[28,12] = 0.9992924598665804
$data[28,13] = 0.9992924598665804
$data[28,14] = 0.9991823959347382
$data[28,15] = 0.9996855382140203
$data[28,16] = 0.9996855382140203
$data[28,17] = 0.9998820773877402
$data[28,18] = 0.9998820773877402
$data[28,19] = 0.9999580722897855
$data[29,0] = 29
$data[29,1] = 'Michael-SNHex'
$data[29,2] = 1.000672536343824
$data[29,3] = 1.000672536343824
$data[29,4] = 0.9985204121046662
$data[29,5] = 1.001793441584459
$data[29,6] = 0.9973098656866343
$data[29,7] = 1.000672536343824
$data[29,8] = 1.000672536343824
$data[29,9] = 0.9973098656866343
$data[29,10] = 1.000672536343824
$data[29,11] = 1.000672536343824
$data[29,12] = 0.9989912010152291
$data[29,13] = 0.9989912010152291
$data[29,14] = 0.9988342713783748
$data[29,15] = 0.9995516461247608
$data[29,16] = 0.9995516461247608
$data[29,17] = 0.9998318686795267
$data[29,18] = 0.9998318686795267
$data[29,19] = 0.9999402214012053

# Write the full A2:T31 block (header row stays at row 1; rows 2-31 are the 30 data rows)
$ws.Range("A2:T31").Value = $data

# Column A carries the bold/centered/bordered "index" style (style used by A2:A29 already);
# copy it onto the two newly-appended rows so A30:A31 match the rest of the column.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false